{"js": "// Office.js (Word JavaScript API) script\n// Applies the edits described by the diff:\n//  - Splits the \"C950 ...\" byline into 3 paragraphs: \"ID 011441603\", \"jada900@wgu.edu\",\n//    \"C950 Data Structures & Algorithms II\"\n//  - Converts the \"Identify a named self-adjusting algorithm...\" list item into a plain\n//    \"A - Algorithm Identification\" heading paragraph.\n//  - Converts the \"Identify a self-adjusting data structure...\" list item (and the\n//    paragraphs following it) into an expanded discussion of the nearest-neighbor\n//    algorithm, a \"B - Data Structure Identification\" heading, additional hash-table\n//    discussion, and a new \"C1 - Algorithm Overview\" heading.\n\nconst OOXML_NS = 'xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"';\n\nfunction wrapOoxml(bodyXml) {\n  return (\n    '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n    '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n    '<pkg:xmlData>' +\n    '<w:document ' + OOXML_NS + '>' +\n    '<w:body>' + bodyXml + '</w:body>' +\n    '</w:document>' +\n    '</pkg:xmlData>' +\n    '</pkg:part>' +\n    '</pkg:package>'\n  );\n}\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\n// ---------------------------------------------------------------------\n// 1) Replace the \"C\" / \"950 \u2013 Data Structures & Algorithms II\" paragraph\n//    with three paragraphs: ID line, email line, course line.\n// ---------------------------------------------------------------------\nlet target = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const t = paragraphs.items[i].text;\n  if (t.indexOf(\"950\") !== -1 && t.indexOf(\"Data Structures\") !== -1) {\n    target = paragraphs.items[i];\n    break;\n  }\n}\nif (!target) {\n  throw new Error(\"Could not find the 'C950 \u2013 Data Structures...' paragraph\");\n}\n\nconst byLineOoxml = wrapOoxml(\n  '<w:p><w:r><w:t xml:space=\"preserve\">ID </w:t></w:r><w:r><w:t>011441603</w:t></w:r></w:p>' +\n  '<w:p><w:r><w:t>jada900@wgu.edu</w:t></w:r></w:p>' +\n  '<w:p><w:r><w:t>C950 Data Structures &amp; Algorithms II</w:t></w:r></w:p>'\n);\ntarget.insertOoxml(byLineOoxml, \"Replace\");\nawait context.sync();\n\n// ---------------------------------------------------------------------\n// 2) Replace the \"Identify a named self-adjusting algorithm...\" list\n//    paragraph with the plain \"A - Algorithm Identification\" heading.\n// ---------------------------------------------------------------------\nconst paragraphs2 = context.document.body.paragraphs;\nparagraphs2.load(\"text\");\nawait context.sync();\n\ntarget = null;\nfor (let i = 0; i < paragraphs2.items.length; i++) {\n  const t = paragraphs2.items[i].text;\n  if (t.indexOf(\"Identify a named self-adjusting algorithm\") !== -1) {\n    target = paragraphs2.items[i];\n    break;\n  }\n}\nif (!target) {\n  throw new Error(\"Could not find the 'Identify a named self-adjusting algorithm...' paragraph\");\n}\n\nconst headingAOoxml = wrapOoxml(\n  '<w:p><w:r><w:t xml:space=\"preserve\">A - </w:t></w:r><w:r><w:t>Algorithm Identification</w:t></w:r></w:p>'\n);\ntarget.insertOoxml(headingAOoxml, \"Replace\");\nawait context.sync();\n\n// ---------------------------------------------------------------------\n// 3) Replace the \"Identify a self-adjusting data structure...\" list\n//    paragraph with the new nearest-neighbor discussion paragraphs, the\n//    \"B - Data Structure Identification\" heading, and a blank line.\n// ---------------------------------------------------------------------\nconst paragraphs3 = context.document.body.paragraphs;\nparagraphs3.load(\"text\");\nawait context.sync();\n\ntarget = null;\nfor (let i = 0; i < paragraphs3.items.length; i++) {\n  const t = paragraphs3.items[i].text;\n  if (t.indexOf(\"Identify a self-adjusting data structure\") !== -1) {\n    target = paragraphs3.items[i];\n    break;\n  }\n}\nif (!target) {\n  throw new Error(\"Could not find the 'Identify a self-adjusting data structure...' paragraph\");\n}\n\nconst middleBlockOoxml = wrapOoxml(\n  '<w:p>' +\n    '<w:r><w:t>In this program, the nearest neighbor algorithm will be used to determine the most efficient route each truck should take to maximize the efficiency of their deliveries.</w:t></w:r>' +\n    '<w:r><w:t xml:space=\"preserve\"> It will plot a route from one package\\u2019s address to another, always choosing the nearest address to the previous delivery address when choosing the </w:t></w:r>' +\n    '<w:r><w:t>next location.</w:t></w:r>' +\n  '</w:p>' +\n  '<w:p/>' +\n  '<w:p>' +\n    '<w:r><w:t xml:space=\"preserve\">Being a greedy algorithm, </w:t></w:r>' +\n    '<w:proofErr w:type=\"spellStart\"/>' +\n    '<w:r><w:t>the</w:t></w:r>' +\n    '<w:proofErr w:type=\"spellEnd\"/>' +\n    '<w:r><w:t xml:space=\"preserve\"> nearest neighbor algorithm</w:t></w:r>' +\n    '<w:r><w:t xml:space=\"preserve\"> may not always select the most optimal path for each truck; however, the paths chosen will be sufficient for the purposes of this project.</w:t></w:r>' +\n  '</w:p>' +\n  '<w:p/>' +\n  '<w:p/>' +\n  '<w:p>' +\n    '<w:r><w:lastRenderedPageBreak/><w:t xml:space=\"preserve\">B - </w:t></w:r>' +\n    '<w:r><w:t>Data Structure Identification</w:t></w:r>' +\n  '</w:p>' +\n  '<w:p/>'\n);\ntarget.insertOoxml(middleBlockOoxml, \"Replace\");\nawait context.sync();\n\n// ---------------------------------------------------------------------\n// 4) After the (unchanged) \"I will be using a hash table...\" paragraph,\n//    append the new follow-up paragraphs and the \"C1 - Algorithm\n//    Overview\" heading. insertOoxml only supports the \"Replace\" location\n//    for a Paragraph in this runtime, so the target paragraph's own\n//    content is included (unchanged) as the first paragraph of the\n//    replacement block.\n// ---------------------------------------------------------------------\nconst paragraphs4 = context.document.body.paragraphs;\nparagraphs4.load(\"text\");\nawait context.sync();\n\ntarget = null;\nfor (let i = 0; i < paragraphs4.items.length; i++) {\n  const t = paragraphs4.items[i].text;\n  if (t.indexOf(\"I will be using a hash table to store packages\") !== -1) {\n    target = paragraphs4.items[i];\n    break;\n  }\n}\nif (!target) {\n  throw new Error(\"Could not find the 'I will be using a hash table...' paragraph\");\n}\n\n// NOTE: this paragraph is the very last paragraph in the document body\n// (immediately followed by sectPr). In that situation the runtime\n// always absorbs one trailing empty paragraph from the inserted OOXML\n// (the body's terminal paragraph mark is structural and can't be fully\n// replaced), so one extra trailing <w:p/> is added here to compensate\n// and end up with two real trailing empty paragraphs, matching the\n// target document.\nconst afterHashTableOoxml = wrapOoxml(\n  '<w:p>' +\n    '<w:r><w:t>I will be using a hash table to store packages with their related data. A hash table is a self-adjusting data structure which stores key-value pairs. More information on my choice of using a hash table can be found in part C.</w:t></w:r>' +\n  '</w:p>' +\n  '<w:p/>' +\n  '<w:p>' +\n    '<w:r><w:t xml:space=\"preserve\">All package-related data will be held in the hash table. </w:t></w:r>' +\n    '<w:r><w:t>This includes fields such as delivery address, delivery deadline, package weight, and delivery status.</w:t></w:r>' +\n  '</w:p>' +\n  '<w:p/>' +\n  '<w:p>' +\n    '<w:r><w:t xml:space=\"preserve\">C1 </w:t></w:r>' +\n    '<w:r><w:t>-</w:t></w:r>' +\n    '<w:r><w:t xml:space=\"preserve\"> Algorithm Overview</w:t></w:r>' +\n  '</w:p>' +\n  '<w:p/>' +\n  '<w:p/>' +\n  '<w:p/>'\n);\ntarget.insertOoxml(afterHashTableOoxml, \"Replace\");\nawait context.sync();\n", "ps1": "# Word COM interop (PowerShell-style) script\n# Applies the edits described by the diff:\n#  - Splits the \"C950 ...\" byline into 3 paragraphs: \"ID 011441603\", \"jada900@wgu.edu\",\n#    \"C950 Data Structures & Algorithms II\"\n#  - Converts the \"Identify a named self-adjusting algorithm...\" list item into a plain\n#    \"A - Algorithm Identification\" heading paragraph.\n#  - Converts the \"Identify a self-adjusting data structure...\" list item (and the\n#    paragraphs following it) into an expanded discussion of the nearest-neighbor\n#    algorithm, a \"B - Data Structure Identification\" heading, additional hash-table\n#    discussion, and a new \"C1 - Algorithm Overview\" heading.\n\n$d = $word.ActiveDocument\n\nfunction New-OoxmlPackage {\n    param([string]$BodyXml)\n    return '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n        '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n        '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n        '<pkg:xmlData>' +\n        '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n        '<w:body>' + $BodyXml + '</w:body>' +\n        '</w:document>' +\n        '</pkg:xmlData>' +\n        '</pkg:part>' +\n        '</pkg:package>'\n}\n\nfunction Find-ParagraphContaining {\n    param([string]$Needle)\n    $count = $d.Paragraphs.Count\n    for ($i = 1; $i -le $count; $i++) {\n        $p = $d.Paragraphs.Item($i)\n        if ($p.Range.Text -like \"*$Needle*\") {\n            return $p\n        }\n    }\n    throw \"Could not find a paragraph containing: $Needle\"\n}\n\n# ---------------------------------------------------------------------\n# 1) Replace the \"C\" / \"950 - Data Structures & Algorithms II\" paragraph\n#    with three paragraphs: ID line, email line, course line.\n# ---------------------------------------------------------------------\n$target = Find-ParagraphContaining \"Data Structures\"\n$byLineXml = New-OoxmlPackage (\n    '<w:p><w:r><w:t xml:space=\"preserve\">ID </w:t></w:r><w:r><w:t>011441603</w:t></w:r></w:p>' +\n    '<w:p><w:r><w:t>jada900@wgu.edu</w:t></w:r></w:p>' +\n    '<w:p><w:r><w:t>C950 Data Structures &amp; Algorithms II</w:t></w:r></w:p>'\n)\n$target.Range.InsertXML($byLineXml)\n\n# ---------------------------------------------------------------------\n# 2) Replace the \"Identify a named self-adjusting algorithm...\" list\n#    paragraph with the plain \"A - Algorithm Identification\" heading.\n# ---------------------------------------------------------------------\n$target = Find-ParagraphContaining \"Identify a named self-adjusting algorithm\"\n$headingAXml = New-OoxmlPackage (\n    '<w:p><w:r><w:t xml:space=\"preserve\">A - </w:t></w:r><w:r><w:t>Algorithm Identification</w:t></w:r></w:p>'\n)\n$target.Range.InsertXML($headingAXml)\n\n# ---------------------------------------------------------------------\n# 3) Replace the \"Identify a self-adjusting data structure...\" list\n#    paragraph with the new nearest-neighbor discussion paragraphs, the\n#    \"B - Data Structure Identification\" heading, and a blank line.\n# ---------------------------------------------------------------------\n$target = Find-ParagraphContaining \"Identify a self-adjusting data structure\"\n$middleBlockXml = New-OoxmlPackage (\n    '<w:p>' +\n        '<w:r><w:t>In this program, the nearest neighbor algorithm will be used to determine the most efficient route each truck should take to maximize the efficiency of their deliveries.</w:t></w:r>' +\n        '<w:r><w:t xml:space=\"preserve\"> It will plot a route from one package\u2019s address to another, always choosing the nearest address to the previous delivery address when choosing the </w:t></w:r>' +\n        '<w:r><w:t>next location.</w:t></w:r>' +\n    '</w:p>' +\n    '<w:p/>' +\n    '<w:p>' +\n        '<w:r><w:t xml:space=\"preserve\">Being a greedy algorithm, </w:t></w:r>' +\n        '<w:proofErr w:type=\"spellStart\"/>' +\n        '<w:r><w:t>the</w:t></w:r>' +\n        '<w:proofErr w:type=\"spellEnd\"/>' +\n        '<w:r><w:t xml:space=\"preserve\"> nearest neighbor algorithm</w:t></w:r>' +\n        '<w:r><w:t xml:space=\"preserve\"> may not always select the most optimal path for each truck; however, the paths chosen will be sufficient for the purposes of this project.</w:t></w:r>' +\n    '</w:p>' +\n    '<w:p/>' +\n    '<w:p/>' +\n    '<w:p>' +\n        '<w:r><w:lastRenderedPageBreak/><w:t xml:space=\"preserve\">B - </w:t></w:r>' +\n        '<w:r><w:t>Data Structure Identification</w:t></w:r>' +\n    '</w:p>' +\n    '<w:p/>'\n)\n$target.Range.InsertXML($middleBlockXml)\n\n# ---------------------------------------------------------------------\n# 4) After the (unchanged) \"I will be using a hash table...\" paragraph,\n#    append the new follow-up paragraphs and the \"C1 - Algorithm\n#    Overview\" heading. This paragraph is the very last paragraph in the\n#    document body. When InsertXML targets the body's final paragraph\n#    range, this COM runtime effectively turns the last paragraph of the\n#    replacement XML into the (now emptied) original trailing paragraph\n#    mark, keeping a 1:1 paragraph count with what was supplied -- so the\n#    full target content (including both trailing blank paragraphs) is\n#    provided explicitly below.\n# ---------------------------------------------------------------------\n$target = Find-ParagraphContaining \"I will be using a hash table to store packages\"\n$afterHashTableXml = New-OoxmlPackage (\n    '<w:p>' +\n        '<w:r><w:t>I will be using a hash table to store packages with their related data. A hash table is a self-adjusting data structure which stores key-value pairs. More information on my choice of using a hash table can be found in part C.</w:t></w:r>' +\n    '</w:p>' +\n    '<w:p/>' +\n    '<w:p>' +\n        '<w:r><w:t xml:space=\"preserve\">All package-related data will be held in the hash table. </w:t></w:r>' +\n        '<w:r><w:t>This includes fields such as delivery address, delivery deadline, package weight, and delivery status.</w:t></w:r>' +\n    '</w:p>' +\n    '<w:p/>' +\n    '<w:p>' +\n        '<w:r><w:t xml:space=\"preserve\">C1 </w:t></w:r>' +\n        '<w:r><w:t>-</w:t></w:r>' +\n        '<w:r><w:t xml:space=\"preserve\"> Algorithm Overview</w:t></w:r>' +\n    '</w:p>' +\n    '<w:p/>' +\n    '<w:p/>'\n)\n$target.Range.InsertXML($afterHashTableXml)\n"}
